# [siva]Add: Report generation for the FaceBook Automation
#
# - Correct the stored password value in Sheet1!B2 ("login@123" -> "Login@123")
# - Update the active selection on Sheet1 to B2 (was G10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the capitalisation of the password test-data value.
$ws.Range("B2").Value = "Login@123"

# Move/save the active selection to B2.
$ws.Activate()
$ws.Range("B2").Select()
